$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the used range
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Swap columns C and D (group-name <-> group-code) for every row, including the header
for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)
    $cValue = $cCell.Value()
    $dValue = $dCell.Value()
    $cCell.Value = $dValue
    $dCell.Value = $cValue
}
